# Update the "cryptos" worksheet with a refreshed GitHub Actions price snapshot.
# Numeric-looking D-column prices are prefixed with a leading apostrophe so
# Excel keeps storing them as text (matching the original inlineStr cells)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.505.52"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.337.91"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'304.47"
$ws.Range("D6").Value = "'101.53"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.515"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "'35.16"
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").Value = "'51.76"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("D12").Value = "'0.0796"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "2.697.48"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "'15.62"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "2.330.27"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "'0.807"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "43.423.66"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'11.77"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "'6.12"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "'67.94"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "'238.41"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("D26").Value = "'2.53"
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'25.00"
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D31").Value = "'9.26"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").Value = "'165.02"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'5.07"
$ws.Range("E34").Value = "  -4.19%  "
$ws.Range("E35").Value = "  -4.92%  "
$ws.Range("E36").Value = "  -4.89%  "
$ws.Range("D37").Value = "'16.93"
$ws.Range("E37").Value = "  -7.54%  "
$ws.Range("D38").Value = "'0.0707"
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("E39").Value = "  -7.09%  "
$ws.Range("E40").Value = "  -6.30%  "
$ws.Range("E41").Value = "  -3.32%  "
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "'2.41"
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("D44").Value = "1.982.32"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").Value = "'18.57"
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("D47").Value = "'2.95"
$ws.Range("E47").Value = "  -7.00%  "
$ws.Range("D48").Value = "'9.93"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("D49").Value = "'56.02"
$ws.Range("E49").Value = "  -4.11%  "
$ws.Range("D50").Value = "'4.91"
$ws.Range("E50").Value = "  +4.21%  "
$ws.Range("D51").Value = "2.563.34"
$ws.Range("E51").Value = "  +0.32%  "

# Row 29/30: swap coin entries (Toncoin <-> InjectiveProtocol), keep rank (column A) fixed
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'34.59"
$ws.Range("E29").Value = "  -6.74%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.07"
$ws.Range("E30").Value = "  -2.13%  "
